# =================================================================
# Commit: feat: add 2022-Q3 data
# - Insert a new "2022-Q3" quarter row into the "总计" summary sheet
# - Insert a new "2022-Q3" worksheet (fund holdings) right after "总计"
# =================================================================

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) "总计" sheet: insert a row for the new quarter at the top of the
#    data block (row 2) and push everything else down by one row.
# -----------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Carry over the bold/centered "index" style from the row below
# (which used to be row 2, now shifted to row 3) onto the new A2.
$wsTotal.Cells.Item(3,1).Copy()
$wsTotal.Cells.Item(2,1).PasteSpecial(-4122)
$wsTotal.Application.CutCopyMode = $false

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 8
$wsTotal.Cells.Item(2,4).Value = 1.07

# -----------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" (so it
#    takes the place the old "2022-Q2" tab used to occupy).
# -----------------------------------------------------------------
$wsAfterAnchor = $wb.Worksheets.Item("总计")
$wsNew = $wb.Worksheets.Add($null, $wsAfterAnchor)
$wsNew.Name = "2022-Q3"

# Re-fetch by name (index-based refs go stale across structural edits).
$ws2 = $wb.Worksheets.Item("2022-Q3")
$wsRef = $wb.Worksheets.Item("2022-Q2")

# Bring over the 9-row range (header + 8 fund rows) from the sibling
# "2022-Q2" sheet so the new sheet starts with matching layout/
# column widths, then we overwrite every cell with the Q3 data and
# fix up the formatting explicitly (cross-sheet paste here does not
# carry cell styles in this host).
$wsRef.Range("A1:H9").Copy()
$ws2.Range("A1").PasteSpecial(-4104)
$ws2.Application.CutCopyMode = $false

$ws2.Cells.Item(1,2).NumberFormat = "@"
$ws2.Cells.Item(1,2).Value = "基金代码"
$ws2.Cells.Item(1,3).NumberFormat = "@"
$ws2.Cells.Item(1,3).Value = "基金名称"
$ws2.Cells.Item(1,4).NumberFormat = "@"
$ws2.Cells.Item(1,4).Value = "基金规模"
$ws2.Cells.Item(1,5).NumberFormat = "@"
$ws2.Cells.Item(1,5).Value = "股票总仓位"
$ws2.Cells.Item(1,6).NumberFormat = "@"
$ws2.Cells.Item(1,6).Value = "仓位占比"
$ws2.Cells.Item(1,7).NumberFormat = "@"
$ws2.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws2.Cells.Item(1,8).NumberFormat = "@"
$ws2.Cells.Item(1,8).Value = "仓位排名"

$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).NumberFormat = "@"
$ws2.Cells.Item(2,2).Value = "007202"
$ws2.Cells.Item(2,3).NumberFormat = "@"
$ws2.Cells.Item(2,3).Value = "天弘优质成长企业精选混合A"
$ws2.Cells.Item(2,4).NumberFormat = "@"
$ws2.Cells.Item(2,4).Value = "5.80"
$ws2.Cells.Item(2,5).NumberFormat = "@"
$ws2.Cells.Item(2,5).Value = "93.00"
$ws2.Cells.Item(2,6).NumberFormat = "@"
$ws2.Cells.Item(2,6).Value = "6.31"
$ws2.Cells.Item(2,7).NumberFormat = "@"
$ws2.Cells.Item(2,7).Value = "0.3660"
$ws2.Cells.Item(2,8).Value = 4

$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,2).NumberFormat = "@"
$ws2.Cells.Item(3,2).Value = "420001"
$ws2.Cells.Item(3,3).NumberFormat = "@"
$ws2.Cells.Item(3,3).Value = "天弘精选混合A"
$ws2.Cells.Item(3,4).NumberFormat = "@"
$ws2.Cells.Item(3,4).Value = "5.23"
$ws2.Cells.Item(3,5).NumberFormat = "@"
$ws2.Cells.Item(3,5).Value = "73.52"
$ws2.Cells.Item(3,6).NumberFormat = "@"
$ws2.Cells.Item(3,6).Value = "4.84"
$ws2.Cells.Item(3,7).NumberFormat = "@"
$ws2.Cells.Item(3,7).Value = "0.2531"
$ws2.Cells.Item(3,8).Value = 4

$ws2.Cells.Item(4,1).Value = 2
$ws2.Cells.Item(4,2).NumberFormat = "@"
$ws2.Cells.Item(4,2).Value = "420005"
$ws2.Cells.Item(4,3).NumberFormat = "@"
$ws2.Cells.Item(4,3).Value = "天弘周期策略混合A"
$ws2.Cells.Item(4,4).NumberFormat = "@"
$ws2.Cells.Item(4,4).Value = "3.10"
$ws2.Cells.Item(4,5).NumberFormat = "@"
$ws2.Cells.Item(4,5).Value = "93.08"
$ws2.Cells.Item(4,6).NumberFormat = "@"
$ws2.Cells.Item(4,6).Value = "8.13"
$ws2.Cells.Item(4,7).NumberFormat = "@"
$ws2.Cells.Item(4,7).Value = "0.2520"
$ws2.Cells.Item(4,8).Value = 3

$ws2.Cells.Item(5,1).Value = 3
$ws2.Cells.Item(5,2).NumberFormat = "@"
$ws2.Cells.Item(5,2).Value = "011851"
$ws2.Cells.Item(5,3).NumberFormat = "@"
$ws2.Cells.Item(5,3).Value = "天弘先进制造混合A"
$ws2.Cells.Item(5,4).NumberFormat = "@"
$ws2.Cells.Item(5,4).Value = "1.92"
$ws2.Cells.Item(5,5).NumberFormat = "@"
$ws2.Cells.Item(5,5).Value = "91.41"
$ws2.Cells.Item(5,6).NumberFormat = "@"
$ws2.Cells.Item(5,6).Value = "5.57"
$ws2.Cells.Item(5,7).NumberFormat = "@"
$ws2.Cells.Item(5,7).Value = "0.1069"
$ws2.Cells.Item(5,8).Value = 5

$ws2.Cells.Item(6,1).Value = 4
$ws2.Cells.Item(6,2).NumberFormat = "@"
$ws2.Cells.Item(6,2).Value = "015458"
$ws2.Cells.Item(6,3).NumberFormat = "@"
$ws2.Cells.Item(6,3).Value = "天弘周期策略混合C"
$ws2.Cells.Item(6,4).NumberFormat = "@"
$ws2.Cells.Item(6,4).Value = "0.68"
$ws2.Cells.Item(6,5).NumberFormat = "@"
$ws2.Cells.Item(6,5).Value = "93.08"
$ws2.Cells.Item(6,6).NumberFormat = "@"
$ws2.Cells.Item(6,6).Value = "8.13"
$ws2.Cells.Item(6,7).NumberFormat = "@"
$ws2.Cells.Item(6,7).Value = "0.0553"
$ws2.Cells.Item(6,8).Value = 3

$ws2.Cells.Item(7,1).Value = 5
$ws2.Cells.Item(7,2).NumberFormat = "@"
$ws2.Cells.Item(7,2).Value = "011852"
$ws2.Cells.Item(7,3).NumberFormat = "@"
$ws2.Cells.Item(7,3).Value = "天弘先进制造混合C"
$ws2.Cells.Item(7,4).NumberFormat = "@"
$ws2.Cells.Item(7,4).Value = "0.67"
$ws2.Cells.Item(7,5).NumberFormat = "@"
$ws2.Cells.Item(7,5).Value = "91.41"
$ws2.Cells.Item(7,6).NumberFormat = "@"
$ws2.Cells.Item(7,6).Value = "5.57"
$ws2.Cells.Item(7,7).NumberFormat = "@"
$ws2.Cells.Item(7,7).Value = "0.0373"
$ws2.Cells.Item(7,8).Value = 5

$ws2.Cells.Item(8,1).Value = 6
$ws2.Cells.Item(8,2).NumberFormat = "@"
$ws2.Cells.Item(8,2).Value = "015460"
$ws2.Cells.Item(8,3).NumberFormat = "@"
$ws2.Cells.Item(8,3).Value = "天弘优质成长企业精选混合C"
$ws2.Cells.Item(8,4).NumberFormat = "@"
$ws2.Cells.Item(8,4).Value = "0.06"
$ws2.Cells.Item(8,5).NumberFormat = "@"
$ws2.Cells.Item(8,5).Value = "93.00"
$ws2.Cells.Item(8,6).NumberFormat = "@"
$ws2.Cells.Item(8,6).Value = "6.31"
$ws2.Cells.Item(8,7).NumberFormat = "@"
$ws2.Cells.Item(8,7).Value = "0.0038"
$ws2.Cells.Item(8,8).Value = 4

$ws2.Cells.Item(9,1).Value = 7
$ws2.Cells.Item(9,2).NumberFormat = "@"
$ws2.Cells.Item(9,2).Value = "015459"
$ws2.Cells.Item(9,3).NumberFormat = "@"
$ws2.Cells.Item(9,3).Value = "天弘精选混合C"
$ws2.Cells.Item(9,4).NumberFormat = "@"
$ws2.Cells.Item(9,4).Value = "0.00"
$ws2.Cells.Item(9,5).NumberFormat = "@"
$ws2.Cells.Item(9,5).Value = "73.52"
$ws2.Cells.Item(9,6).NumberFormat = "@"
$ws2.Cells.Item(9,6).Value = "4.84"
$ws2.Cells.Item(9,7).Value = 0
$ws2.Cells.Item(9,8).Value = 4

# -----------------------------------------------------------------
# 3) Formatting: bold + centered + top-aligned + thin border for the
#    header row (B1:H1) and the index column (A1:A9), matching the
#    style used by every other quarter sheet in the workbook.
# -----------------------------------------------------------------
$hdr = $ws2.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$idxCol = $ws2.Range("A1:A9")
$idxCol.Font.Bold = $true
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
$idxCol.Borders.LineStyle = 1

Write-Host "2022-Q3 quarter added."
